$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$np = $s.NotesPage
$sh = $np.Shapes.Item(2)
$sh.TextFrame.TextRange.Text = "Overview of PMU from slides. `rThis is the HPC used in our project to help analyze the performance of event schedular."
